$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was updated from 45189
# (2023-09-20) to 45190 (2023-09-21) for every data row (rows 2-90).
for ($r = 2; $r -le 90; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value = 45190
    }
}
